# Insert a new price-record row at row 516 (a weekly Acelga quote from
# Femacal de La Calera), pushing the existing rows 516:585 down to 517:586.
# This mirrors the commit's "Fruta / hortaliza, semanal" weekly refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 516 onward down by one to make room for the new record.
$ws.Rows.Item(516).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(516, 1).Value  = 3
$ws.Cells.Item(516, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(516, 3).Value  = "Coquimbo"
$ws.Cells.Item(516, 4).Value  = 45127
$ws.Cells.Item(516, 5).Value  = 5
$ws.Cells.Item(516, 6).Value  = 100112009
$ws.Cells.Item(516, 7).Value  = "Acelga"
$ws.Cells.Item(516, 8).Value  = "Sin especificar"
$ws.Cells.Item(516, 9).Value  = "Primera"
$ws.Cells.Item(516, 10).Value = 210
$ws.Cells.Item(516, 11).Value = 3000
$ws.Cells.Item(516, 12).Value = 3300
$ws.Cells.Item(516, 13).Value = 3143
$ws.Cells.Item(516, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(516, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(516, 16).Value = 524
$ws.Cells.Item(516, 17).Value = 6
$ws.Cells.Item(516, 18).Value = "Hortaliza"
